{"js": "// The document's Title paragraph and Abstract paragraph were each split\n// across several runs (one run per word plus separate single-space runs).\n// The edit collapses each of those paragraphs down to a single run whose\n// text is the same, fully concatenated sentence.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text,style\");\nawait context.sync();\n\nlet titlePara = null;\nlet abstractPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n    const para = paragraphs.items[i];\n    if (titlePara === null && para.style === \"Title\") {\n        titlePara = para;\n    } else if (abstractPara === null && para.style === \"Abstract\") {\n        abstractPara = para;\n    }\n}\n\nif (titlePara !== null) {\n    titlePara.insertText(\"Answers: Trigonometry (degrees)\", Word.InsertLocation.replace);\n}\n\nif (abstractPara !== null) {\n    abstractPara.insertText(\n        \"Answers to the questions on trigonometry, using degrees to measure angles.\",\n        Word.InsertLocation.replace\n    );\n}\n\nawait context.sync();\n", "ps1": "# The document's Title paragraph and Abstract paragraph were each split\n# across several runs (one run per word plus separate single-space runs).\n# This collapses each of those paragraphs down to a single run whose text\n# is the same, fully concatenated sentence.\n#\n# Find/Replace (scoped to each paragraph's own Range) is used instead of a\n# plain Range.Text assignment because it reliably rewrites the whole match\n# as one run even when the replacement text is identical to the text\n# already there (which is the case here - only the run layout changes).\n\n$d = $word.ActiveDocument\n\n$newTitleText = \"Answers: Trigonometry (degrees)\"\n$newAbstractText = \"Answers to the questions on trigonometry, using degrees to measure angles.\"\n\n$titlePara = $null\n$abstractPara = $null\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Style.NameLocal\n    if ($titlePara -eq $null -and $styleName -eq \"Title\") {\n        $titlePara = $p\n    } elseif ($abstractPara -eq $null -and $styleName -eq \"Abstract\") {\n        $abstractPara = $p\n    }\n}\n\nif ($titlePara -ne $null) {\n    $titleRange = $titlePara.Range\n    $titleFind = $titleRange.Find\n    [void]$titleFind.Execute($newTitleText, $false, $false, $false, $false, $false, $true, 1, $false, $newTitleText, 2)\n}\n\nif ($abstractPara -ne $null) {\n    $abstractRange = $abstractPara.Range\n    $abstractFind = $abstractRange.Find\n    [void]$abstractFind.Execute($newAbstractText, $false, $false, $false, $false, $false, $true, 1, $false, $newAbstractText, 2)\n}\n"}
